# Daily attendance processing - 2026-01-10 21:55:40
# Swap the order of recorders in the "Recorded By" column (G) wherever a
# session was recorded jointly by the instructor and the automated
# "System" actor: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$used.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com") | Out-Null
